$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$row = 86

$rng = $ws.Range("A" + $row + ":J" + $row)
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-05-26"
$ws.Cells.Item($row, 2).Value = "35.5"
$ws.Cells.Item($row, 3).Value = "35.4"
$ws.Cells.Item($row, 4).Value = "0.94"
$ws.Cells.Item($row, 5).Value = "0.258"
$ws.Cells.Item($row, 6).Value = "0.09"
$ws.Cells.Item($row, 7).Value = "5,373"
$ws.Cells.Item($row, 8).Value = "8,045"
$ws.Cells.Item($row, 9).Value = "8,095"
$ws.Cells.Item($row, 10).Value = "7.2241"

# Values are now stored as text; drop the temporary "@" number format so
# the new cells end up with the default (unstyled) appearance, matching
# the rest of the sheet.
$rng.ClearFormats()
